$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels and data values for columns K, L, M (rows 1-26)
$Kvals = @("ARANCELES CON IVA", "202109,28", "451263,69", "499314,65", "579279,91", "138602,48", "130733,60", "135633,12", "135694,17", "332815,19", "129315,81", "208566,12", "107259,25", "208566,12", "88874,24", "366067,18", "571076,12", "331362,80", "331362,80", "258752,11", "214752,42", "188379,33", "118293,00", "309593,32", "133051,21", "258752,11")
$Lvals = @("IVA ARANCELES", "167032,46", "372945,20", "412656,74", "478743,73", "114547,50", "108044,30", "112093,49", "112143,94", "275053,88", "106872,57", "172368,69", "88644,01", "172368,69", "73449,79", "302534,86", "471963,74", "273853,55", "273853,55", "213844,72", "177481,34", "155685,40", "97762,81", "255862,25", "109959,68", "213844,72")
$Mvals = @("ARANCELES SIN IVA", "35076,82", "78318,49", "86657,91", "100536,18", "24054,98", "22689,30", "23539,63", "23550,23", "57761,31", "22443,24", "36197,43", "18615,24", "36197,43", "15424,45", "63532,32", "99112,38", "57509,25", "57509,25", "44907,39", "37271,08", "32693,93", "20530,19", "53731,07", "23091,53", "44907,39")

# Copy header formatting (bold, centered, thin border) from existing header cell J1
# onto the new header cells before writing their text.
$ws.Range("J1").Copy()
$ws.Range("K1:M1").PasteSpecial(-4122)

for ($i = 0; $i -lt 26; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 11).Value = $Kvals[$i]
    $ws.Cells.Item($r, 12).Value = $Lvals[$i]
    $ws.Cells.Item($r, 13).Value = $Mvals[$i]
}
